$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.8708544531477855
$arr[0,1] = 0.2115725435158708
$arr[0,2] = 0.2831738088168265
$arr[0,3] = 0
$arr[0,4] = 1.169114993605803
$arr[0,5] = 0.002429833495797397
$arr[0,6] = 0
$arr[0,7] = 0.4923759621168529
$arr[0,8] = 0.2654150915616924
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.4317726924522987
$arr[0,12] = 0
$arr[0,13] = 2.48527645876618
$ws.Range("B2:O2").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.7692415225463947
$arr[0,1] = 0.1848225412675788
$arr[0,2] = 0.2801412553420448
$arr[0,3] = 0
$arr[0,4] = 1.174571900294438
$arr[0,5] = 0.002432651108401205
$arr[0,6] = 0
$arr[0,7] = 0.5016695270689357
$arr[0,8] = 0.2667356277561979
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.4014935418317975
$arr[0,12] = 0
$arr[0,13] = 2.505528179058203
$ws.Range("B3:O3").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.7067111210932069
$arr[0,1] = 0.1683397309201098
$arr[0,2] = 0.2783788234778228
$arr[0,3] = 0
$arr[0,4] = 1.178754154028084
$arr[0,5] = 0.002434473377209276
$arr[0,6] = 0
$arr[0,7] = 0.5077817493025414
$arr[0,8] = 0.2677064589397489
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.3829769276011916
$arr[0,12] = 0
$arr[0,13] = 2.519955727439338
$ws.Range("B4:O4").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.68119594481783
$arr[0,1] = 0.161608608636044
$arr[0,2] = 0.2776857411989369
$arr[0,3] = 0
$arr[0,4] = 1.180667365264647
$arr[0,5] = 0.002435239231489994
$arr[0,6] = 0
$arr[0,7] = 0.5103744091218747
$arr[0,8] = 0.2681422856184739
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.3754505206266856
$arr[0,12] = 0
$arr[0,13] = 2.526335325404546
$ws.Range("B5:O5").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.6769571910970456
$arr[0,1] = 0.1604900622611183
$arr[0,2] = 0.2775721753159672
$arr[0,3] = 0
$arr[0,4] = 1.18099766347482
$arr[0,5] = 0.002435367808081179
$arr[0,6] = 0
$arr[0,7] = 0.5108110658308682
$arr[0,8] = 0.2682170816497518
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.3742019426270602
$arr[0,12] = 0
$arr[0,13] = 2.5274248411706
$ws.Range("B6:O6").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.7063671481509743
$arr[0,1] = 0.1682490097992968
$arr[0,2] = 0.2783693744916746
$arr[0,3] = 0
$arr[0,4] = 1.178779110702862
$arr[0,5] = 0.002434483611624096
$arr[0,6] = 0
$arr[0,7] = 0.5078163025743869
$arr[0,8] = 0.2677121739050676
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.382875345276851
$arr[0,12] = 0
$arr[0,13] = 2.520039740588388
$ws.Range("B7:O7").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 0.835848215888916
$arr[0,1] = 0.2023614483340737
$arr[0,2] = 0.2821075688690371
$arr[0,3] = 0
$arr[0,4] = 1.17082377930317
$arr[0,5] = 0.002430785905451211
$arr[0,6] = 0
$arr[0,7] = 0.4954960144037024
$arr[0,8] = 0.2658371811195366
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.4213171403526275
$arr[0,12] = 0
$arr[0,13] = 2.491845093043281
$ws.Range("B8:O8").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.088595119354409
$arr[0,1] = 0.268780106976152
$arr[0,2] = 0.2902250376943698
$arr[0,3] = 0
$arr[0,4] = 1.16183384608135
$arr[0,5] = 0.00242426344211579
$arr[0,6] = 0
$arr[0,7] = 0.4745654034996427
$arr[0,8] = 0.2634316527636358
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.4972809644886738
$arr[0,12] = 0
$arr[0,13] = 2.452409466049204
$ws.Range("B9:O9").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.273520081180948
$arr[0,1] = 0.3172739881642599
$arr[0,2] = 0.2966653295670341
$arr[0,3] = 0
$arr[0,4] = 1.159276201641688
$arr[0,5] = 0.002419911175705021
$arr[0,6] = 0
$arr[0,7] = 0.4611677476810492
$arr[0,8] = 0.2624419511899418
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5534309965348854
$arr[0,12] = 0
$arr[0,13] = 2.433160262274328
$ws.Range("B10:O10").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.357469503390575
$arr[0,1] = 0.3392665222817186
$arr[0,2] = 0.2996979891482283
$arr[0,3] = 0
$arr[0,4] = 1.158995319932899
$arr[0,5] = 0.002418025762527864
$arr[0,6] = 0
$arr[0,7] = 0.4555051696288501
$arr[0,8] = 0.2621611507552402
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5790462620778385
$arr[0,12] = 0
$arr[0,13] = 2.426527646690545
$ws.Range("B11:O11").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.389232657518335
$arr[0,1] = 0.3475844858517689
$arr[0,2] = 0.3008611090876059
$arr[0,3] = 0
$arr[0,4] = 1.159016156517396
$arr[0,5] = 0.002417325316357145
$arr[0,6] = 0
$arr[0,7] = 0.4534232425436677
$arr[0,8] = 0.2620792243361052
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5887561589899946
$arr[0,12] = 0
$arr[0,13] = 2.424322474816023
$ws.Range("B12:O12").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.382393106512836
$arr[0,1] = 0.3457935207435412
$arr[0,2] = 0.3006099571766754
$arr[0,3] = 0
$arr[0,4] = 1.159006007221109
$arr[0,5] = 0.002417475569597468
$arr[0,6] = 0
$arr[0,7] = 0.45386884519964
$arr[0,8] = 0.2620957824609107
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5866645234756618
$arr[0,12] = 0
$arr[0,13] = 2.424783751891511
$ws.Range("B13:O13").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.360083221204491
$arr[0,1] = 0.3399510515676525
$arr[0,2] = 0.2997933852481651
$arr[0,3] = 0
$arr[0,4] = 1.158994483269396
$arr[0,5] = 0.002417967865873242
$arr[0,6] = 0
$arr[0,7] = 0.4553326369335515
$arr[0,8] = 0.2621539212212696
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5798449042566602
$arr[0,12] = 0
$arr[0,13] = 2.426340078129215
$ws.Range("B14:O14").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.346414249798102
$arr[0,1] = 0.3363710388039181
$arr[0,2] = 0.2992951255145755
$arr[0,3] = 0
$arr[0,4] = 1.159003998139923
$arr[0,5] = 0.002418271170182426
$arr[0,6] = 0
$arr[0,7] = 0.4562373808130609
$arr[0,8] = 0.2621927126653247
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5756689732956346
$arr[0,12] = 0
$arr[0,13] = 2.427333313122489
$ws.Range("B15:O15").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.268030147370951
$arr[0,1] = 0.3158353290214677
$arr[0,2] = 0.2964692022536894
$arr[0,3] = 0
$arr[0,4] = 1.159312342391232
$arr[0,5] = 0.002420036287188468
$arr[0,6] = 0
$arr[0,7] = 0.4615465252105473
$arr[0,8] = 0.2624637151003029
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5517583957123975
$arr[0,12] = 0
$arr[0,13] = 2.433636547057631
$ws.Range("B16:O16").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.219898315426633
$arr[0,1] = 0.3032197323246919
$arr[0,2] = 0.2947618943754264
$arr[0,3] = 0
$arr[0,4] = 1.15972775939953
$arr[0,5] = 0.002421143276009029
$arr[0,6] = 0
$arr[0,7] = 0.4649143476711686
$arr[0,8] = 0.2626733902701019
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.537108236248713
$arr[0,12] = 0
$arr[0,13] = 2.43804810765414
$ws.Range("B17:O17").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.192197891924934
$arr[0,1] = 0.2959572387719618
$arr[0,2] = 0.2937895901477106
$arr[0,3] = 0
$arr[0,4] = 1.160049754258964
$arr[0,5] = 0.002421788881017906
$arr[0,6] = 0
$arr[0,7] = 0.466892085699282
$arr[0,8] = 0.2628099338275263
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5286886871300425
$arr[0,12] = 0
$arr[0,13] = 2.440785374610584
$ws.Range("B18:O18").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.18281626171796
$arr[0,1] = 0.2934972077181897
$arr[0,2] = 0.2934620522064364
$arr[0,3] = 0
$arr[0,4] = 1.160173032830826
$arr[0,5] = 0.002422009001174903
$arr[0,6] = 0
$arr[0,7] = 0.467568689281368
$arr[0,8] = 0.2628589021701657
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5258391601666688
$arr[0,12] = 0
$arr[0,13] = 2.441746459688574
$ws.Range("B19:O19").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.225023726789061
$arr[0,1] = 0.3045633428999963
$arr[0,2] = 0.2949426375267876
$arr[0,3] = 0
$arr[0,4] = 1.159674939303088
$arr[0,5] = 0.002421024515039814
$arr[0,6] = 0
$arr[0,7] = 0.4645516278088522
$arr[0,8] = 0.2626494195313143
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5386670665938453
$arr[0,12] = 0
$arr[0,13] = 2.437557797842771
$ws.Range("B20:O20").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.366636912352988
$arr[0,1] = 0.3416674046472963
$arr[0,2] = 0.3000328334609748
$arr[0,3] = 0
$arr[0,4] = 1.158994413701905
$arr[0,5] = 0.002417822900141103
$arr[0,6] = 0
$arr[0,7] = 0.4549009912859638
$arr[0,8] = 0.2621361817168761
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5818477252109204
$arr[0,12] = 0
$arr[0,13] = 2.425874622036559
$ws.Range("B21:O21").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.459033098557825
$arr[0,1] = 0.3658578083650355
$arr[0,2] = 0.303445316951013
$arr[0,3] = 0
$arr[0,4] = 1.159291202531904
$arr[0,5] = 0.002415809228532983
$arr[0,6] = 0
$arr[0,7] = 0.4489573817638863
$arr[0,8] = 0.2619430320891425
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.6101265820291104
$arr[0,12] = 0
$arr[0,13] = 2.420025727315192
$ws.Range("B22:O22").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.409734355997159
$arr[0,1] = 0.3529524915738307
$arr[0,2] = 0.3016161924371943
$arr[0,3] = 0
$arr[0,4] = 1.159064855542297
$arr[0,5] = 0.002416876776444567
$arr[0,6] = 0
$arr[0,7] = 0.4520962473299548
$arr[0,8] = 0.2620330863382208
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5950284829533388
$arr[0,12] = 0
$arr[0,13] = 2.422983565002937
$ws.Range("B23:O23").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.222706617872291
$arr[0,1] = 0.3039559264915397
$arr[0,2] = 0.294860894729652
$arr[0,3] = 0
$arr[0,4] = 1.159698560230709
$arr[0,5] = 0.00242107817837478
$arr[0,6] = 0
$arr[0,7] = 0.46471548410479
$arr[0,8] = 0.2626602068716792
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.5379623099471473
$arr[0,12] = 0
$arr[0,13] = 2.437778840870379
$ws.Range("B24:O24").Value = $arr

$arr = New-Object 'object[,]' 1,14
$arr[0,0] = 1.020351163992018
$arr[0,1] = 0.2508644810438341
$arr[0,2] = 0.2879451507513551
$arr[0,3] = 0
$arr[0,4] = 1.163556369144843
$arr[0,5] = 0.0024259503891661
$arr[0,6] = 0
$arr[0,7] = 0.479880784606614
$arr[0,8] = 0.2639460384589185
$arr[0,9] = 0
$arr[0,10] = 0
$arr[0,11] = 0.4766701415866237
$arr[0,12] = 0
$arr[0,13] = 2.461374444022283
$ws.Range("B25:O25").Value = $arr

Write-Output "done"